$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3:D3").Value = 5
$ws.Range("C5:D5").Value = 980
$ws.Range("C7:D7").Value = 263
$ws.Range("C9:D9").Value = 216
$ws.Range("C11:D11").Value = 1010
$ws.Range("C13:D13").Value = 471
$ws.Range("C16:D16").Value = 47
$ws.Range("C17:D17").Value = 573
$ws.Range("C19:D19").Value = 262
$ws.Range("C21:D21").Value = 195
$ws.Range("C23:D23").Value = 782
$ws.Range("C25:D25").Value = 9
$ws.Range("C27:D27").Value = 324
$ws.Range("C29:D29").Value = 907
$ws.Range("C32:D32").Value = 96
$ws.Range("C33:D33").Value = 215
$ws.Range("C35:D35").Value = 77
$ws.Range("C37:D37").Value = 145
$ws.Range("C39:D39").Value = 2720
$ws.Range("C41:D41").Value = 735
$ws.Range("C42:D42").Value = 653
$ws.Range("C44:D44").Value = 60
$ws.Range("C46:D46").Value = 2064
$ws.Range("C48:D48").Value = 186
$ws.Range("C50:D50").Value = 1457
$ws.Range("C52:D52").Value = 148
$ws.Range("C54:D54").Value = 212
$ws.Range("C56:D56").Value = 1004
$ws.Range("C58:D58").Value = 4369
$ws.Range("C60:D60").Value = 57
$ws.Range("C62:D62").Value = 4220
$ws.Range("C64:D64").Value = 530
$ws.Range("C66:D66").Value = 196
$ws.Range("C68:D68").Value = 97
$ws.Range("C70:D70").Value = 1643
$ws.Range("C72:D72").Value = 149
$ws.Range("C74:D74").Value = 203
$ws.Range("C76:D76").Value = 115
$ws.Range("C78:D78").Value = 1391
$ws.Range("C80:D80").Value = 211
$ws.Range("C82:D82").Value = 111
$ws.Range("C84:D84").Value = 106
$ws.Range("C86:D86").Value = 99
$ws.Range("C88:D88").Value = 283
$ws.Range("C90:D90").Value = 709
$ws.Range("C92:D92").Value = 1402
$ws.Range("C94:D94").Value = 736
$ws.Range("C96:D96").Value = 107
$ws.Range("C98:D98").Value = 622
$ws.Range("C100:D100").Value = 183
$ws.Range("C102:D102").Value = 332
$ws.Range("C104:D104").Value = 12
$ws.Range("C106:D106").Value = 356
$ws.Range("C108:D108").Value = 101
$ws.Range("C110:D110").Value = 116
$ws.Range("C112:D112").Value = 213
$ws.Range("C114:D114").Value = 362
$ws.Range("C116:D116").Value = 4896
$ws.Range("C118:D118").Value = 710
$ws.Range("C120:D120").Value = 501
$ws.Range("C122:D122").Value = 162
$ws.Range("C124:D124").Value = 425
$ws.Range("C126:D126").Value = 104
$ws.Range("C128:D128").Value = 105
$ws.Range("C130:D130").Value = 22
$ws.Range("C132:D132").Value = 102
$ws.Range("C134:D134").Value = 328
$ws.Range("C136:D136").Value = 2522
$ws.Range("C138:D138").Value = 164
$ws.Range("C140:D140").Value = 542
$ws.Range("C142:D142").Value = 1812
$ws.Range("C144:D144").Value = 1813
$ws.Range("C145").Value = 681.3194444444445